# Scheduled runner: refresh Leve profit-calc sheets with latest market-board prices.
# (currentAveragePrice / NQ / HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H-N)
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 133.63158
$ws.Range("I33").Value = 142.73334
$ws.Range("K33").Value = 142.73334
$ws.Range("M33").Value = 86.26666

$ws.Range("H74").Value = 4241.9165
$ws.Range("I74").Value = 3380.6
$ws.Range("K74").Value = 3380.6
$ws.Range("M74").Value = -2444.6

$ws.Range("H77").Value = 4241.9165
$ws.Range("I77").Value = 3380.6
$ws.Range("K77").Value = 16903
$ws.Range("M77").Value = -12223

$ws.Range("H98").Value = 3363
$ws.Range("I98").Value = 2919.4
$ws.Range("K98").Value = 2919.4
$ws.Range("M98").Value = -1421.4

$ws.Range("H103").Value = 1474.75
$ws.Range("I103").Value = 1400
$ws.Range("K103").Value = 4200
$ws.Range("M103").Value = -3614

$ws.Range("H122").Value = 3363
$ws.Range("I122").Value = 2919.4
$ws.Range("K122").Value = 8758.200000000001
$ws.Range("M122").Value = -6308.200000000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 4500
$ws.Range("I63").Value = 4500
$ws.Range("J63").Value = 4500
$ws.Range("K63").Value = 4500
$ws.Range("L63").Value = 4500
$ws.Range("M63").Value = -3814
$ws.Range("N63").Value = -5872

$ws.Range("H66").Value = 4500
$ws.Range("I66").Value = 4500
$ws.Range("J66").Value = 4500
$ws.Range("K66").Value = 22500
$ws.Range("L66").Value = 22500
$ws.Range("M66").Value = -19068
$ws.Range("N66").Value = -29364

$ws.Range("H122").Value = 4091.0386
$ws.Range("I122").Value = 2378.8333
$ws.Range("J122").Value = 7943.5
$ws.Range("K122").Value = 7136.499899999999
$ws.Range("L122").Value = 23830.5
$ws.Range("M122").Value = -4686.499899999999
$ws.Range("N122").Value = -28730.5

$ws.Range("H139").Value = 68381.75
$ws.Range("J139").Value = 68381.75
$ws.Range("L139").Value = 68381.75
$ws.Range("N139").Value = -78661.75

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1617.4166
$ws.Range("I107").Value = 2173.3635
$ws.Range("J107").Value = 1147
$ws.Range("K107").Value = 2173.3635
$ws.Range("L107").Value = 1147
$ws.Range("M107").Value = -253.3634999999999
$ws.Range("N107").Value = -4987

$ws.Range("H123").Value = 59699.5
$ws.Range("J123").Value = 59699.5
$ws.Range("L123").Value = 59699.5
$ws.Range("N123").Value = -69499.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H18").Value = 55000
$ws.Range("J18").Value = 55000
$ws.Range("L18").Value = 55000
$ws.Range("N18").Value = -55460

$ws.Range("H99").Value = 17981.092
$ws.Range("I99").Value = 24891.8
$ws.Range("J99").Value = 12222.167
$ws.Range("K99").Value = 24891.8
$ws.Range("L99").Value = 12222.167
$ws.Range("M99").Value = -23393.8
$ws.Range("N99").Value = -15218.167

$ws.Range("H117").Value = 41500
$ws.Range("J117").Value = 41500
$ws.Range("L117").Value = 41500
$ws.Range("N117").Value = -50678

$ws.Range("H122").Value = 4932
$ws.Range("I122").Value = 3133.4
$ws.Range("J122").Value = 5931.222
$ws.Range("K122").Value = 9400.200000000001
$ws.Range("L122").Value = 17793.666
$ws.Range("M122").Value = -6950.200000000001
$ws.Range("N122").Value = -22693.666

$ws.Range("H123").Value = 46977.8
$ws.Range("J123").Value = 52472.25
$ws.Range("L123").Value = 52472.25
$ws.Range("N123").Value = -62272.25

$ws.Range("H126").Value = 17981.092
$ws.Range("I126").Value = 24891.8
$ws.Range("J126").Value = 12222.167
$ws.Range("K126").Value = 74675.39999999999
$ws.Range("L126").Value = 36666.501
$ws.Range("M126").Value = -72205.39999999999
$ws.Range("N126").Value = -41606.501

$ws.Range("H129").Value = 99592.336
$ws.Range("J129").Value = 99592.336
$ws.Range("L129").Value = 99592.336
$ws.Range("N129").Value = -109592.336

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 10954
$ws.Range("I3").Value = 3931
$ws.Range("K3").Value = 11793
$ws.Range("M3").Value = -11681

$ws.Range("H122").Value = 486.46667
$ws.Range("J122").Value = 483.91666
$ws.Range("L122").Value = 4355.24994
$ws.Range("N122").Value = -9255.24994

$ws.Range("H129").Value = 1912.05
$ws.Range("I129").Value = 1000.5
$ws.Range("J129").Value = 2823.6
$ws.Range("K129").Value = 3001.5
$ws.Range("L129").Value = 8470.799999999999
$ws.Range("M129").Value = 1998.5
$ws.Range("N129").Value = -18470.8

$ws.Range("H131").Value = 7144705
$ws.Range("J131").Value = 5850052
$ws.Range("L131").Value = 17550156
$ws.Range("N131").Value = -17560236

$ws.Range("H134").Value = 20004.666
$ws.Range("I134").Value = 23339
$ws.Range("J134").Value = 3333
$ws.Range("K134").Value = 70017
$ws.Range("L134").Value = 9999
$ws.Range("M134").Value = -64947
$ws.Range("N134").Value = -20139

$ws.Range("H137").Value = 7870344.5
$ws.Range("I137").Value = 4241.6665
$ws.Range("K137").Value = 12724.9995
$ws.Range("M137").Value = -7624.999500000002

$ws.Range("H138").Value = 66360.875
$ws.Range("I138").Value = 146650
$ws.Range("J138").Value = 3913.7778
$ws.Range("K138").Value = 439950
$ws.Range("L138").Value = 11741.3334
$ws.Range("M138").Value = -434810
$ws.Range("N138").Value = -22021.3334

$ws.Range("H139").Value = 773017.3
$ws.Range("I139").Value = 1001972.7
$ws.Range("K139").Value = 3005918.1
$ws.Range("M139").Value = -3000778.1

$ws.Range("H141").Value = 6373.75
$ws.Range("I141").Value = 5319.1055
$ws.Range("J141").Value = 10381.4
$ws.Range("K141").Value = 15957.3165
$ws.Range("L141").Value = 31144.2
$ws.Range("M141").Value = -10777.3165
$ws.Range("N141").Value = -41504.2

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 7150.567
$ws.Range("I102").Value = 6864.2085
$ws.Range("J102").Value = 8296
$ws.Range("K102").Value = 6864.2085
$ws.Range("L102").Value = 8296
$ws.Range("M102").Value = -5242.2085
$ws.Range("N102").Value = -11540

$ws.Range("H122").Value = 788283.9
$ws.Range("I122").Value = 848690.3
$ws.Range("K122").Value = 2546070.9
$ws.Range("M122").Value = -2543620.9

$ws.Range("H138").Value = 59949
$ws.Range("J138").Value = 59949
$ws.Range("L138").Value = 59949
$ws.Range("N138").Value = -70229

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4995.0977
$ws.Range("I7").Value = 3828.5667
$ws.Range("J7").Value = 8176.5454
$ws.Range("K7").Value = 3828.5667
$ws.Range("L7").Value = 8176.5454
$ws.Range("M7").Value = -3716.5667
$ws.Range("N7").Value = -8400.545399999999

$ws.Range("H46").Value = 6263.08
$ws.Range("I46").Value = 3413.2856
$ws.Range("K46").Value = 3413.2856
$ws.Range("M46").Value = -3225.2856

$ws.Range("H126").Value = 4995.0977
$ws.Range("I126").Value = 3828.5667
$ws.Range("J126").Value = 8176.5454
$ws.Range("K126").Value = 11485.7001
$ws.Range("L126").Value = 24529.6362
$ws.Range("M126").Value = -9015.7001
$ws.Range("N126").Value = -29469.6362

$ws.Range("H132").Value = 6062.5654
$ws.Range("I132").Value = 3731.2856
$ws.Range("K132").Value = 11193.8568
$ws.Range("M132").Value = -8663.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 2503.4255
$ws.Range("J122").Value = 3542.7144
$ws.Range("L122").Value = 10628.1432
$ws.Range("N122").Value = -15528.1432

$ws.Range("H136").Value = 7567.05
$ws.Range("I136").Value = 3228.5652
$ws.Range("J136").Value = 8862.960999999999
$ws.Range("K136").Value = 9685.695599999999
$ws.Range("L136").Value = 26588.883
$ws.Range("M136").Value = -7135.695599999999
$ws.Range("N136").Value = -31688.883
